$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "f085c14f-2fc9-4f17-a2b6-6ddaa2dafb7b"
$ws.Range("B5").Value = "kerchief"
$ws.Range("C5").Value = "Handicrafts"
$ws.Range("D5").Value = "Handcrafted kerchief made with traditional techniques"
$ws.Range("E5").Value = "உயர்தர kerchief, பாரம்பரிய கைவினைத் திறனால் உருவாக்கப்பட்டது"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 890
$ws.Range("H5").Value = "2025-07-16 11:19:01"

$ws.Range("A6").Value = "e74a9b1c-2c57-46af-a6a9-0f32747e47fa"
$ws.Range("B6").Value = "கம்பர் கட்"
$ws.Range("C6").Value = "Handicrafts"
$ws.Range("D6").Value = "Handcrafted கம்பர் கட் made with traditional techniques"
$ws.Range("E6").Value = "உயர்தர கம்பர் கட், பாரம்பரிய கைவினைத் திறனால் உருவாக்கப்பட்டது"
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 22
$ws.Range("H6").Value = "2025-07-16 11:31:01"
